# Weekly fruit/vegetable price update: a new record for Cilantro at
# "Vega Modelo de Temuco" was added to the top of the (date-descending)
# data block, pushing the previous rows 248-256 down to 249-257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 248, shifting rows 248:256
# down to 249:257 (same as Excel's Insert > Entire Row).
$ws.Rows.Item(248).Insert()

# Populate the newly inserted row 248 with the new weekly record.
$ws.Range("A248").Value = 10
$ws.Range("B248").Value = "Vega Modelo de Temuco"
$ws.Range("C248").Value = "La Araucanía"
$ws.Range("D248").Value = 44509
$ws.Range("E248").Value = 9
$ws.Range("F248").Value = 100112040
$ws.Range("G248").Value = "Cilantro"
$ws.Range("H248").Value = "Sin especificar"
$ws.Range("I248").Value = "Primera"
$ws.Range("J248").Value = 30
$ws.Range("K248").Value = 4000
$ws.Range("L248").Value = 4000
$ws.Range("M248").Value = 4000
$ws.Range("N248").Value = "$/docena de atados (2 kilos)"
$ws.Range("O248").Value = "Provincia de Cautín"
$ws.Range("P248").Value = 2000
$ws.Range("Q248").Value = 2
$ws.Range("R248").Value = "Hortaliza"
